# Sofya issues by mail on 2023-08-03
#
# Adds three new key/value translation pairs ("study", "copy" /
# "copy of this print", "plate_at" / "plate at") to the `translations`
# sheet, reusing the existing A (key) / D (English) column layout, and
# updates the sheet's active selection to reflect where the editor was
# last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# New rows 44-46: column A holds the translation key, column D the
# English text shown in the UI (mirrors every other row in the sheet).
$ws.Range("A44").Value = "study"
$ws.Range("D44").Value = "study"

$ws.Range("A45").Value = "copy"
$ws.Range("D45").Value = "copy of this print"

$ws.Range("A46").Value = "plate_at"
$ws.Range("D46").Value = "plate at"

# Move the active selection to D47, just below the rows that were added.
$ws.Range("D47").Select() | Out-Null
